$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking strings
# (e.g. "1.016") are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '27.383.54'
$ws.Range('E2').Value = '  +1.90%  '
$ws.Range('D3').Value = '1.843.27'
$ws.Range('E3').Value = '  +1.49%  '
$ws.Range('D4').Value = '1.016'
$ws.Range('E4').Value = '  +1.47%  '
$ws.Range('D5').Value = '316.82'
$ws.Range('E5').Value = '  +2.41%  '
$ws.Range('D6').Value = '1.014'
$ws.Range('E6').Value = '  +1.24%  '
$ws.Range('D7').Value = '0.4739'
$ws.Range('E7').Value = '  +1.22%  '
$ws.Range('D8').Value = '0.3703'
$ws.Range('E8').Value = '  +0.40%  '
$ws.Range('D9').Value = '0.07457'
$ws.Range('E9').Value = '  +1.16%  '
$ws.Range('D10').Value = '0.8872'
$ws.Range('E10').Value = '  +1.95%  '
$ws.Range('D11').Value = '20.51'
$ws.Range('E11').Value = '  +0.47%  '
$ws.Range('D12').Value = '1.877.20'
$ws.Range('E12').Value = '  -0.31%  '
$ws.Range('D13').Value = '0.07378'
$ws.Range('E13').Value = '  +4.35%  '
$ws.Range('D14').Value = '5.494'
$ws.Range('E14').Value = '  +2.40%  '
$ws.Range('D15').Value = '93.38'
$ws.Range('E15').Value = '  +1.10%  '
$ws.Range('D16').Value = '6.590'
$ws.Range('E16').Value = '  +1.52%  '
$ws.Range('D17').Value = '1.017'
$ws.Range('D18').Value = '0.000008859'
$ws.Range('E18').Value = '  +1.78%  '
$ws.Range('E19').Value = '  +1.26%  '
$ws.Range('E20').Value = '  +0.80%  '
$ws.Range('D21').Value = '27.408.25'
$ws.Range('E21').Value = '  +1.72%  '
$ws.Range('D22').Value = '5.342'
$ws.Range('E22').Value = '  -0.07%  '
$ws.Range('D23').Value = '10.73'
$ws.Range('E23').Value = '  +1.61%  '
$ws.Range('D24').Value = '2.069.82'
$ws.Range('E24').Value = '  -0.66%  '
$ws.Range('D25').Value = '1.911'
$ws.Range('E25').Value = '  +0.47%  '
$ws.Range('D26').Value = '152.50'
$ws.Range('E26').Value = '  +0.78%  '
$ws.Range('D27').Value = '18.69'
$ws.Range('E27').Value = '  +1.83%  '
$ws.Range('D28').Value = '2.177'
$ws.Range('E28').Value = '  +0.40%  '
$ws.Range('D29').Value = '5.289'
$ws.Range('E29').Value = '  -0.59%  '
$ws.Range('D30').Value = '118.16'
$ws.Range('E30').Value = '  +2.21%  '
$ws.Range('D31').Value = '0.08977'
$ws.Range('E31').Value = '  +0.49%  '
$ws.Range('D32').Value = '0.7621'
$ws.Range('E32').Value = '  -0.62%  '
$ws.Range('D33').Value = '1.177'
$ws.Range('E33').Value = '  +1.36%  '
$ws.Range('E34').Value = '  +1.45%  '
$ws.Range('D35').Value = '2.953'
$ws.Range('E35').Value = '  +1.79%  '
$ws.Range('D36').Value = '1.015'
$ws.Range('E36').Value = '  +1.35%  '
$ws.Range('E37').Value = '  +1.98%  '
$ws.Range('D38').Value = '0.05368'
$ws.Range('D39').Value = '0.01971'
$ws.Range('E39').Value = '  +0.37%  '
$ws.Range('D40').Value = '3.009'
$ws.Range('E40').Value = '  +2.22%  '
$ws.Range('D41').Value = '7.337'
$ws.Range('D42').Value = '2.400'
$ws.Range('E42').Value = '  +2.18%  '
$ws.Range('D43').Value = '0.5367'
$ws.Range('E43').Value = '  +0.94%  '
$ws.Range('D44').Value = '0.1669'
$ws.Range('E44').Value = '  +0.21%  '
$ws.Range('E45').Value = '  +1.62%  '
$ws.Range('D46').Value = '0.4967'
$ws.Range('E46').Value = '  +0.67%  '
$ws.Range('D47').Value = '10.50'
$ws.Range('E47').Value = '  +0.32%  '
$ws.Range('D48').Value = '1.015'
$ws.Range('E48').Value = '  +1.36%  '
$ws.Range('D49').Value = '104.81'
$ws.Range('E49').Value = '  +0.86%  '
$ws.Range('D50').Value = '1.683'
$ws.Range('E50').Value = '  +0.86%  '
$ws.Range('D51').Value = '0.06336'
$ws.Range('E51').Value = '  +0.95%  '

# Restore the original (default) cell style on column D so the
# workbook XML does not pick up a stray style reference.
$ws.Range("D2:D51").Style = "Normal"

